# Update patient admission data on the active sheet before migrating to v2 system.
# Values that could otherwise be auto-interpreted by Excel as a number/date
# (plain digit strings, ISO dates) are entered with a leading apostrophe so
# they are stored as literal text, matching the original workbook's text
# (shared-string) cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient name and clinical record number
$ws.Range("A6").Value = "MORALES  CAMEY  ROCIO  JAZMIN"
$ws.Range("G6").Value = "'9211"

# Date of birth, age, place of birth
$ws.Range("A9").Value = "'2007-05-07"
$ws.Range("D9").Value = "'10"
$ws.Range("E9").Value = "GUATEMALA"

# Occupation, nationality, identification document
$ws.Range("C11").Value = "ESTUDIA"
$ws.Range("E11").Value = "GUATEMALTECA"
$ws.Range("G11").Value = "NO PRESENTO"

# Emergency contact info: name, relationship, address, phone
$ws.Range("A13").Value = "DORIBEL CAMEY"
$ws.Range("D13").Value = "MADRE"
$ws.Range("E13").Value = "LOTE 15 SANTA MARTA "
$ws.Range("G13").Value = "'47201313"

# Area of urgency and assistance date
$ws.Range("E14").Value = "Area de urgencia: CIRUGIA"
$ws.Range("A15").Value = "'2017-10-19"
